$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure columns D and E keep their existing text format so that
# numeric-looking values (e.g. "0.572", "1.975.40") are stored as
# literal text rather than being auto-converted to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '59.863.16'
$ws.Range('E2').Value = '  +1.27%  '
$ws.Range('D3').Value = '2.650.34'
$ws.Range('E3').Value = '  +2.31%  '
$ws.Range('D5').Value = '536.88'
$ws.Range('E5').Value = '  +1.47%  '
$ws.Range('D6').Value = '145.37'
$ws.Range('E6').Value = '  +3.98%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('D8').Value = '0.572'
$ws.Range('E8').Value = '  +1.32%  '
$ws.Range('D9').Value = '2.668.19'
$ws.Range('E9').Value = '  +2.49%  '
$ws.Range('E10').Value = '  +3.92%  '
$ws.Range('E11').Value = '  +1.89%  '
$ws.Range('D12').Value = '0.338'
$ws.Range('E12').Value = '  +1.51%  '
$ws.Range('E13').Value = '  -1.31%  '
$ws.Range('D14').Value = '3.127.52'
$ws.Range('E14').Value = '  +2.53%  '
$ws.Range('D15').Value = '59.786.21'
$ws.Range('E15').Value = '  +1.24%  '
$ws.Range('D16').Value = '21.22'
$ws.Range('E16').Value = '  +3.83%  '
$ws.Range('D17').Value = '2.664.92'
$ws.Range('E17').Value = '  +3.60%  '
$ws.Range('D18').Value = '0.0000135'
$ws.Range('E18').Value = '  +1.34%  '
$ws.Range('D19').Value = '344.58'
$ws.Range('E19').Value = '  -0.67%  '
$ws.Range('D20').Value = '4.42'
$ws.Range('E20').Value = '  +2.03%  '
$ws.Range('D21').Value = '10.28'
$ws.Range('E21').Value = '  +1.93%  '
$ws.Range('D22').Value = '6.36'
$ws.Range('E22').Value = '  -0.53%  '
$ws.Range('E23').Value = '  +0.07%  '
$ws.Range('D24').Value = '67.02'
$ws.Range('E24').Value = '  -0.71%  '
$ws.Range('D25').Value = '0.415'
$ws.Range('E25').Value = '  +2.40%  '
$ws.Range('D26').Value = '0.166'
$ws.Range('E26').Value = '  -0.29%  '
$ws.Range('E27').Value = '  +0.50%  '
$ws.Range('D28').Value = '7.34'
$ws.Range('E28').Value = '  +2.65%  '
$ws.Range('D29').Value = '0.0₃0752'
$ws.Range('E29').Value = '  +2.41%  '
$ws.Range('E30').Value = '  -0.07%  '
$ws.Range('E31').Value = '  +3.15%  '
$ws.Range('D32').Value = '5.86'
$ws.Range('E32').Value = '  +0.21%  '
$ws.Range('D33').Value = '19.07'
$ws.Range('E33').Value = '  +1.63%  '
$ws.Range('D34').Value = '150.29'
$ws.Range('E34').Value = '  +1.17%  '
$ws.Range('D35').Value = '4.05'
$ws.Range('E35').Value = '  +1.95%  '
$ws.Range('E36').Value = '  +3.50%  '
$ws.Range('D37').Value = '1.47'
$ws.Range('E37').Value = '  +0.09%  '
$ws.Range('D38').Value = '0.845'
$ws.Range('E38').Value = '  +1.96%  '
$ws.Range('D39').Value = '0.823'
$ws.Range('E39').Value = '  -0.46%  '
$ws.Range('D40').Value = '294.58'
$ws.Range('E40').Value = '  +9.49%  '
$ws.Range('D41').Value = '3.60'
$ws.Range('E41').Value = '  +2.15%  '
$ws.Range('E42').Value = '  +0.03%  '
$ws.Range('D43').Value = '0.605'
$ws.Range('E43').Value = '  +1.72%  '
$ws.Range('E44').Value = '  +5.45%  '
$ws.Range('B45').Value = 'WhiteBITCoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D45').Value = '10.74'
$ws.Range('E45').Value = '  -0.05%  '
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').Value = '19.36'
$ws.Range('E46').Value = '  +5.00%  '
$ws.Range('D47').Value = '0.0955'
$ws.Range('E47').Value = '  -0.27%  '
$ws.Range('D48').Value = '0.0227'
$ws.Range('E48').Value = '  +2.72%  '
$ws.Range('D49').Value = '1.975.40'
$ws.Range('E49').Value = '  +1.27%  '
$ws.Range('D50').Value = '18.48'
$ws.Range('E50').Value = '  +1.73%  '
$ws.Range('D51').Value = '4.57'
$ws.Range('E51').Value = '  -0.90%  '
